$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that held only the teacher's name in columns B/C (with no label in
# column A) is removed; everything below it shifts up by one row, carrying
# its row heights along with it.
$ws.Rows(13).Delete()

# After the shift, a handful of cells need content that doesn't simply slide
# up from the row above/below - the underlying data changed. Fix those cells
# explicitly.

# Row 10 (Objetivos:) and row 18 (Método:) both now show the teacher name
# that used to live in the deleted row.
$ws.Range("B10").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C10").Value = "5817650 - Érica Leonor Romão"
$ws.Range("B18").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C18").Value = "5817650 - Érica Leonor Romão"

# Row 13 (Programa resumido:) now reads "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now reads "01/01/2020" - copy the value from the
# existing "Ativação:" cell (row 8) instead of typing it, so it stays plain
# text instead of being reinterpreted as a date value/format.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0
